# JPSDI93_test_2.xlsx - "changes about test outputs"
#
# 1. Update scored output values + summary row on the "outputs" sheet.
# 2. Append a new empty "outputs2" sheet after "outputs1".

$wb = $excel.ActiveWorkbook

# --- 1. Update the "outputs" sheet --------------------------------------
$ws = $wb.Worksheets.Item("outputs")

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 5

$ws.Range("A6").Value = "psdi"
$ws.Range("B6").Value = "level_1,level_4"

# --- 2. Add a new "outputs2" sheet after "outputs1" ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "outputs2"

# Match the look & feel of the other (empty) output sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Restore the original selected/active sheet so adding the new sheet
# doesn't change which tab is active in the saved workbook.
$wb.Worksheets.Item(1).Activate()
